$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.546.47"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "3.452.92"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.96%  "
$ws.Range("D7").Value = "3.454.02"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "4.041.57"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.65%  "
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "3.451.54"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "61.682.66"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +8.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.565"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").Value = "3.591.73"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.69%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("D37").Value = "3.480.73"
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "167.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "28.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.54%  "
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.802"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.53%  "
$ws.Range("E47").Value = "  +4.49%  "
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "2.599.47"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("E51").Value = "  +2.49%  "
